$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers in K1 and L1
$ws.Range("K1").Value = "volumen"
$ws.Range("L1").Value = "flujo"

# Fill rows 2 through 20 with the new data
for ($r = 2; $r -le 20; $r++) {
    $kCell = $ws.Cells.Item($r, 11)  # column K
    $lCell = $ws.Cells.Item($r, 12)  # column L
    $kCell.Value = 200
    $lCell.Value = 12
    $kCell.NumberFormat = "@"
    $lCell.NumberFormat = "@"
}

# Update the active selection to K5, as reflected in the edited workbook
$ws.Range("K5").Select()
